$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaRange = $metaPara.Range

$metaTail = ": Discover the unique gameplay style of Blazing X, an oriental-flavored slot game with a free spin feature and 25x multiplier. Play for free online."

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>' + $metaTail + '</w:t></w:r></w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$metaRange.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Remove the old bold "Play Blazing X Free Slot Game Online" paragraph
#    that used to live at the very end of the document, right before the
#    italic meta-description-like paragraph.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPara = $d.Paragraphs($count - 1)
$boldRange = $d.Range($boldPara.Range.Start, $boldPara.Range.End)
$boldRange.Delete()

# ---------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    feature-image prompt, keeping its italic formatting.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range

# Clear the existing text but keep the paragraph mark in place.
$clearRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$clearRange.Delete()

$imageText = 'Create a feature image for &quot;Blazing X&quot;. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be holding a dragon and standing in front of the Blazing X game grid. The background should be highlighted in fiery orange and red hues to represent the theme of the game. Use creative illustrations and vivid colors to make the image pop and capture the attention of slot players looking for a new and exciting game to play.'

$lastPara2 = $d.Paragraphs($d.Paragraphs.Count)
$lastRange2 = $lastPara2.Range

$imageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>' + $imageText + '</w:t></w:r></w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$lastRange2.InsertXML($imageXml)

# Inserting XML at the very end of the body story leaves a stray empty
# trailing paragraph behind (the split point has nothing to merge into);
# clean it back up so the edited paragraph is once again the last one.
if ($d.Paragraphs.Count -gt $count - 1) {
    $trailing = $d.Paragraphs($d.Paragraphs.Count)
    if ($trailing.Range.Text -eq "") {
        $trailRange = $d.Range($trailing.Range.Start - 1, $trailing.Range.End)
        $trailRange.Delete()
    }
}

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
Write-Output ("Para2: " + $d.Paragraphs(2).Range.Text)
Write-Output ("LastPara: " + $d.Paragraphs($d.Paragraphs.Count).Range.Text)
